$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws.Range("D2").Value = 748.04
$ws.Range("E2").Value = -748.04

$ws.Range("D4").Value = 748.04
$ws.Range("E4").Value = 16751.96
$ws.Range("F4").Value = 0.04274514285714286
